# Update "想去人数" (want-to-go count) figures in column F for the
# "展览" (Exhibition) and "全部类型" (All types) sheets, matching the
# newly scraped totals. Sheets "演出" and "本地生活" are unaffected.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibition) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F4").Value  = 1134
$ws1.Range("F5").Value  = 121
$ws1.Range("F6").Value  = 74
$ws1.Range("F7").Value  = 275
$ws1.Range("F8").Value  = 59
$ws1.Range("F9").Value  = 1164
$ws1.Range("F10").Value = 16399
$ws1.Range("F11").Value = 281
$ws1.Range("F14").Value = 6374
$ws1.Range("F15").Value = 639
$ws1.Range("F16").Value = 127
$ws1.Range("F18").Value = 27
$ws1.Range("F19").Value = 126
$ws1.Range("F20").Value = 1274
$ws1.Range("F21").Value = 42
$ws1.Range("F24").Value = 34
$ws1.Range("F26").Value = 2
$ws1.Range("F27").Value = 17
$ws1.Range("F28").Value = 223
$ws1.Range("F29").Value = 896
$ws1.Range("F31").Value = 5052
$ws1.Range("F32").Value = 503
$ws1.Range("F33").Value = 11342
$ws1.Range("F34").Value = 1250
$ws1.Range("F36").Value = 151
$ws1.Range("F38").Value = 3840

# --- Sheet 4: 全部类型 (All types) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F4").Value  = 1134
$ws4.Range("F5").Value  = 121
$ws4.Range("F6").Value  = 74
$ws4.Range("F7").Value  = 275
$ws4.Range("F8").Value  = 59
$ws4.Range("F9").Value  = 1164
$ws4.Range("F10").Value = 16399
$ws4.Range("F11").Value = 281
$ws4.Range("F14").Value = 6374
$ws4.Range("F15").Value = 639
$ws4.Range("F16").Value = 127
$ws4.Range("F18").Value = 27
$ws4.Range("F19").Value = 126
$ws4.Range("F20").Value = 1274
$ws4.Range("F21").Value = 42
$ws4.Range("F24").Value = 34
$ws4.Range("F26").Value = 2
$ws4.Range("F27").Value = 17
$ws4.Range("F28").Value = 223
$ws4.Range("F29").Value = 896
$ws4.Range("F31").Value = 5052
$ws4.Range("F32").Value = 503
$ws4.Range("F34").Value = 11342
$ws4.Range("F35").Value = 1250
$ws4.Range("F37").Value = 151
$ws4.Range("F39").Value = 3840
